# Updates cryptos list prices and volume(1h) percentages, matching the
# GitHub Actions automated data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = "new price"; E = "new volume text" }
# Only D and/or E keys present when that column changed for the row.
$updates = @{
    2  = @{ D = "37.454.49"; E = "  +0.85%  " }
    3  = @{ D = "2.018.22";  E = "  +0.95%  " }
    5  = @{ D = "260.38";    E = "  +5.93%  " }
    6  = @{ D = "0.616";     E = "  -0.98%  " }
    7  = @{ E = "  -0.07%  " }
    8  = @{ D = "56.58";     E = "  -5.38%  " }
    9  = @{ D = "0.389";     E = "  +1.13%  " }
    10 = @{ D = "0.0776";    E = "  -3.46%  " }
    11 = @{ E = "  -1.33%  " }
    12 = @{ D = "14.37";     E = "  -3.98%  " }
    13 = @{ D = "2.314.23";  E = "  +0.83%  " }
    14 = @{ E = "  -3.73%  " }
    15 = @{ E = "  -6.24%  " }
    16 = @{ E = "  -2.56%  " }
    17 = @{ D = "2.024.08";  E = "  +1.19%  " }
    18 = @{ D = "37.359.59"; E = "  +0.79%  " }
    19 = @{ D = "69.98";     E = "  -0.18%  " }
    20 = @{ D = "0.0₃0840";  E = "  -2.56%  " }
    21 = @{ D = "5.21";      E = "  +0.65%  " }
    22 = @{ D = "228.85";    E = "  -0.51%  " }
    23 = @{ E = "  +8.14%  " }
    25 = @{ E = "  +0.78%  " }
    26 = @{ D = "165.13";    E = "  +0.57%  " }
    27 = @{ D = "9.03";      E = "  -4.19%  " }
    28 = @{ D = "19.79";     E = "  +0.99%  " }
    29 = @{ E = "  -6.42%  " }
    30 = @{ D = "1.32";      E = "  -2.75%  " }
    31 = @{ E = "  -0.58%  " }
    32 = @{ E = "  -2.22%  " }
    33 = @{ D = "0.0650";    E = "  -0.93%  " }
    34 = @{ D = "4.61";      E = "  +3.20%  " }
    35 = @{ E = "  +1.15%  " }
    36 = @{ E = "  +0.95%  " }
    37 = @{ D = "3.38";      E = "  +1.93%  " }
    38 = @{ E = "  -0.13%  " }
    39 = @{ D = "5.27";      E = "  -1.64%  " }
    40 = @{ E = "  +4.04%  " }
    41 = @{ E = "  +2.49%  " }
    42 = @{ E = "  -0.27%  " }
    43 = @{ D = "0.0937";    E = "  -4.52%  " }
    44 = @{ D = "1.402.31";  E = "  +2.55%  " }
    45 = @{ D = "90.39";     E = "  -0.49%  " }
    46 = @{ D = "15.84";     E = "  -4.28%  " }
    47 = @{ E = "  -1.26%  " }
    48 = @{ D = "7.14";      E = "  -2.77%  " }
    49 = @{ E = "  +2.21%  " }
    50 = @{ D = "2.205.16";  E = "  +0.82%  " }
    51 = @{ E = "  -5.45%  " }
}

# Price/volume text must stay plain text (matches the source data's
# inline-string cells), not get auto-coerced into numbers by Excel's
# "looks like a number" detection. Forcing the Text number format before
# the write, then resetting the style back to Normal afterwards, keeps the
# literal string value without leaving a stray number format behind.
foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    if ($cols.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cols["D"]
        $cell.Style = "Normal"
    }
    if ($cols.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cols["E"]
        $cell.Style = "Normal"
    }
}
